# Refined metadata to be additional tab
#
# 1) Refresh the per-row "time_taken" timestamps on the existing "data" sheet
#    (column F, rows 2-14) to reflect the new data pull.
# 2) Add a new "metadata" worksheet (after "data") summarising the panel
#    query that produced the data tab.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1) Update the "time_taken" timestamps on the data sheet -----------------
$newTimes = @{
    2  = "2021-10-05 14:35:46.351613"
    3  = "2021-10-05 14:35:46.351621"
    4  = "2021-10-05 14:35:46.351624"
    5  = "2021-10-05 14:35:46.351627"
    6  = "2021-10-05 14:35:46.351630"
    7  = "2021-10-05 14:35:46.351633"
    8  = "2021-10-05 14:35:46.351635"
    9  = "2021-10-05 14:35:46.351638"
    10 = "2021-10-05 14:35:46.351641"
    11 = "2021-10-05 14:35:46.351643"
    12 = "2021-10-05 14:35:46.351646"
    13 = "2021-10-05 14:35:46.351649"
    14 = "2021-10-05 14:35:46.351651"
}

foreach ($row in $newTimes.Keys) {
    $dataSheet.Range("F$row").Value = $newTimes[$row]
}

# --- 2) Add the "metadata" worksheet right after "data" ----------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the bold/boxed header style from the data sheet (B1:F1) for the new
# header row (B1:G1), and the index-column style (A2) for the new A2 cell.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)   # xlPasteFormats

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

$metaSheet.Application.CutCopyMode = $false

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Skeletal Muscle Channelopathies"
$metaSheet.Range("C2").Value = 302
$metaSheet.Range("E2").Value = "2021-01-17T02:45:09.201879Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:46.347944"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/302/?format=json"

# D2 must hold the literal text "1.0" (not the number 1). Stage it as a text
# formula in a scratch cell, then copy only the resulting value over so the
# target cell ends up as a plain/no-style text cell (matching the diff) and
# no extra number-format styles get introduced.
$metaSheet.Range("ZZ1").Formula = '="1.0"'
$metaSheet.Range("ZZ1").Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)      # xlPasteValues
$metaSheet.Range("ZZ1").Clear()
$metaSheet.Application.CutCopyMode = $false

Write-Host "Inserted 'metadata' worksheet and refreshed data timestamps."
